# Replace exact text occurrences while preserving surrounding run structure
# (e.g. leading empty <w:r/> placeholder runs) by performing the edit as a
# tracked insert + delete and then accepting all revisions. A direct
# Range.Text / Find-replace assignment causes the engine to coalesce
# same-formatted adjacent runs (dropping empty runs); going through the
# track-changes path keeps run boundaries intact.
function Replace-ExactText {
    param($doc, [string]$oldText, [string]$newText)

    $searchRange = $doc.Content
    $found = $searchRange.Find.Execute(
        $oldText, $true, $false, $false, $false, $false,
        $true, 1, $false, $null, 0
    )
    if (-not $found) {
        return $false
    }

    $startPos = $searchRange.Start
    $endPos = $searchRange.End

    $wasTracking = $doc.TrackRevisions
    $doc.TrackRevisions = $true

    # Insert the new text right after the old text first (this keeps any
    # preceding empty run's position untouched), then delete the old text.
    $insPoint = $doc.Range($endPos, $endPos)
    $insPoint.InsertBefore($newText)

    $oldRange = $doc.Range($startPos, $endPos)
    $oldRange.Delete()

    $doc.TrackRevisions = $wasTracking
    $doc.AcceptAllRevisions()

    return $true
}

$d = $word.ActiveDocument

# 1. Title heading + bold recap line (two occurrences of identical text)
Replace-ExactText $d "Play Forbidden Dragon Free: Chinese Mythical Theme Slot" "Play Forbidden Dragon for Free | Review and Gameplay"
Replace-ExactText $d "Play Forbidden Dragon Free: Chinese Mythical Theme Slot" "Play Forbidden Dragon for Free | Review and Gameplay"

# 2. "What we like" bullet list items
Replace-ExactText $d "Chinese theme featuring mythical creatures" "Traditional Chinese symbols on the paytable"
Replace-ExactText $d "Colossal Reel engine for more winning combinations" "Bonus features include wild and scatter symbols"
Replace-ExactText $d "Free spins and Wild Transfer features" "Clear and well-structured user interface"
Replace-ExactText $d "Playable on a range of devices" "Built on the Colossal Reels game engine"

# 3. "What we don't like" bullet list items
Replace-ExactText $d "High volatility may not be ideal for all players" "Free spins bonus can be challenging to trigger"
Replace-ExactText $d "Free spins can be challenging to trigger" "High volatility may not appeal to all players"

# 4. Final italic summary line
Replace-ExactText $d "Get immersed in Chinese mythology as you play the Forbidden Dragon online slot for free. Enjoy the Colossal Reel engine and unique bonus features." "Play Forbidden Dragon for free and explore its Chinese theme and Colossal Reel game engine."
